# Updated cryptos list (Price / Volume(1h) columns) per the latest data pull.
# For cells whose new text looks like a plain number (e.g. "586.60"), we assign
# it through Formula with a leading apostrophe so Excel keeps it as text (matching
# the original inline-string cell type) instead of silently converting it to a
# numeric value, then restore the cell's default (unstyled) look by copying the
# Style from an untouched data cell in the same column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.986.55"
$ws.Range("E2").Value = "  -4.27%  "

$ws.Range("D3").Value = "3.804.94"
$ws.Range("E3").Value = "  -4.93%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Formula = "'586.60"
$ws.Range("D5").Style = $ws.Range("D2").Style
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").Formula = "'163.65"
$ws.Range("D6").Style = $ws.Range("D2").Style
$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("E7").Value = "  -3.74%  "

$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("D9").Formula = "'0.735"
$ws.Range("D9").Style = $ws.Range("D2").Style
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("D10").Formula = "'0.173"
$ws.Range("D10").Style = $ws.Range("D2").Style
$ws.Range("E10").Value = "  +2.41%  "

$ws.Range("D11").Formula = "'52.01"
$ws.Range("D11").Style = $ws.Range("D2").Style
$ws.Range("E11").Value = "  -4.15%  "

$ws.Range("E12").Value = "  -0.80%  "

$ws.Range("D13").Formula = "'11.06"
$ws.Range("D13").Style = $ws.Range("D2").Style
$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("D14").Value = "4.422.05"
$ws.Range("E14").Value = "  -4.50%  "

$ws.Range("D15").Value = "3.828.65"
$ws.Range("E15").Value = "  -4.42%  "

$ws.Range("D16").Formula = "'20.54"
$ws.Range("D16").Style = $ws.Range("D2").Style
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("E17").Value = "  -4.20%  "

$ws.Range("E18").Value = "  -6.91%  "

$ws.Range("E19").Value = "  -2.52%  "

$ws.Range("D20").Value = "69.940.68"
$ws.Range("E20").Value = "  -3.95%  "

$ws.Range("D21").Formula = "'433.74"
$ws.Range("D21").Style = $ws.Range("D2").Style
$ws.Range("E21").Value = "  -0.31%  "

$ws.Range("D22").Formula = "'4.64"
$ws.Range("D22").Style = $ws.Range("D2").Style
$ws.Range("E22").Value = "  -3.47%  "

$ws.Range("D23").Formula = "'92.53"
$ws.Range("D23").Style = $ws.Range("D2").Style
$ws.Range("E23").Value = "  -4.19%  "

$ws.Range("E24").Value = "  -6.73%  "

$ws.Range("D25").Formula = "'13.66"
$ws.Range("D25").Style = $ws.Range("D2").Style
$ws.Range("E25").Value = "  -4.20%  "

$ws.Range("D26").Formula = "'11.03"
$ws.Range("D26").Style = $ws.Range("D2").Style
$ws.Range("E26").Value = "  -2.92%  "

$ws.Range("E27").Value = "  -11.91%  "

$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("E30").Value = "  -5.07%  "

$ws.Range("D31").Formula = "'7.92"
$ws.Range("D31").Style = $ws.Range("D2").Style
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").Formula = "'13.26"
$ws.Range("D32").Style = $ws.Range("D2").Style
$ws.Range("E32").Value = "  -4.13%  "

$ws.Range("D33").Formula = "'47.62"
$ws.Range("D33").Style = $ws.Range("D2").Style
$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("E34").Value = "  -6.36%  "

$ws.Range("D35").Formula = "'68.46"
$ws.Range("D35").Style = $ws.Range("D2").Style
$ws.Range("E35").Value = "  -3.47%  "

$ws.Range("D36").Value = "0.0₃0963"
$ws.Range("E36").Value = "  +8.92%  "

$ws.Range("D37").Formula = "'625.38"
$ws.Range("D37").Style = $ws.Range("D2").Style
$ws.Range("E37").Value = "  -7.05%  "

$ws.Range("D38").Formula = "'0.417"
$ws.Range("D38").Style = $ws.Range("D2").Style
$ws.Range("E38").Value = "  -5.97%  "

$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").Formula = "'1.00"
$ws.Range("D40").Style = $ws.Range("D2").Style
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("E41").Value = "  -2.91%  "

$ws.Range("E42").Value = "  -4.45%  "

$ws.Range("D43").Formula = "'3.13"
$ws.Range("D43").Style = $ws.Range("D2").Style
$ws.Range("E43").Value = "  +17.99%  "

$ws.Range("E44").Value = "  -5.76%  "

$ws.Range("D45").Formula = "'2.75"
$ws.Range("D45").Style = $ws.Range("D2").Style
$ws.Range("E45").Value = "  +4.73%  "

$ws.Range("D46").Formula = "'9.74"
$ws.Range("D46").Style = $ws.Range("D2").Style
$ws.Range("E46").Value = "  -9.55%  "

$ws.Range("E47").Value = "  -5.76%  "

$ws.Range("E48").Value = "  -16.07%  "

$ws.Range("D49").Formula = "'3.26"
$ws.Range("D49").Style = $ws.Range("D2").Style
$ws.Range("E49").Value = "  -4.74%  "

$ws.Range("D50").Value = "2.804.44"
$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("E51").Value = "  -0.45%  "

